$d = $word.ActiveDocument

# --- Simple text replacements (Find & Replace) ---

$d.Content.Find.Execute(
    "Placed at the beginning of each dose", $true, $false, $false, $false,
    $false, $true, 1, $false, "Placed for each dose", 2)

$d.Content.Find.Execute(
    "Removed that day at the end of treatment", $true, $false, $false, $false,
    $false, $true, 1, $false, "Removed that day", 2)

$d.Content.Find.Execute(
    "May shower within 24 hrs", $true, $false, $false, $false,
    $false, $true, 1, $false, "May shower in 24 hrs", 2)

$d.Content.Find.Execute(
    "Placed underneath the skin below the right collarbone", $true, $false, $false, $false,
    $false, $true, 1, $false, "Implanted under the skin", 2)

$d.Content.Find.Execute(
    "Incision in the neck (1/4”)", $true, $false, $false, $false,
    $false, $true, 1, $false, "Neck incision (1/4”)", 2)

$d.Content.Find.Execute(
    "Critical to good communication with your cancer care team", $true, $false, $false, $false,
    $false, $true, 1, $false, "Critical to good communication with your care team", 2)

$d.Content.Find.Execute(
    "Important to reduce the risk of complications from cancer treatment", $true, $false, $false, $false,
    $false, $true, 1, $false, "Reduces risk of complications from treatment", 2)

$d.Content.Find.Execute(
    "Working hard enough that you can’t carry a conversation", $true, $false, $false, $false,
    $false, $true, 1, $false, "Working hard enough that you can’t converse", 2)

$d.Content.Find.Execute(
    "Start slow an build up", $true, $false, $false, $false,
    $false, $true, 1, $false, "Start slowly and build up", 2)

$d.Content.Find.Execute(
    "Smoking makes it more difficult to get through cancer treatment", $true, $false, $false, $false,
    $false, $true, 1, $false, "Smoking makes cancer treatment more difficult", 2)

$d.Content.Find.Execute(
    "American Lung Asssociation fredomfromsmoking.org", $true, $false, $false, $false,
    $false, $true, 1, $false, "American Lung Assn fredomfromsmoking.org", 2)

$d.Content.Find.Execute(
    "1:1 Smoking Cessation Counseling Clinics (Metro Charlotte)", $true, $false, $false, $false,
    $false, $true, 1, $false, "1:1 Smoking Cessation Counseling (Metro Charlotte)", 2)

# --- Remove the empty "FirstParagraph" paragraph (just a manual line break)
# that sits right before the "Critical to good communication..." bullet list,
# under the "My Atrium Patient Portal" heading. ---

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "First Paragraph" -and $p.Range.Text.Trim() -eq "") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -like "Critical to good communication*") {
            $p.Range.Delete()
        }
    }
}
